$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 89; this pushes the existing data (rows 89-127)
# down to rows 90-128, preserving each row's formatting (matches the diff:
# dimension grows from A1:T127 to A1:T128 and every row from 89 downward is
# shifted by one with a brand-new record now occupying row 89).
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new market record.
$ws.Range("A89").Value = 3
$ws.Range("B89").Value = "Femacal de La Calera"
$ws.Range("C89").Value = "Coquimbo"
$ws.Range("D89").Value = 44704
$ws.Range("E89").Value = 5
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100107
$ws.Range("H89").Value = "Otros"
$ws.Range("I89").Value = 100107011
$ws.Range("J89").Value = "Tuna"
$ws.Range("K89").Value = "Sin especificar"
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 68
$ws.Range("N89").Value = 17000
$ws.Range("O89").Value = 17000
$ws.Range("P89").Value = 17000
$ws.Range("Q89").Value = "$/caja 20 kilos"
$ws.Range("R89").Value = "Provincia de Limarí"
$ws.Range("S89").Value = 850
$ws.Range("T89").Value = 20
